# Auto commit update: refresh Metrics input values for the new day,
# letting dependent formulas on the "today" sheet (and TODAY()-1) recalc
# automatically, then restore the original active sheet/tab + selections.

$wb = $excel.ActiveWorkbook

# --- Update the raw metric inputs on the "Metrics" sheet ---------------
$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value  = 387693.93000000005
$wsMetrics.Range("B3").Value  = 331917.64999999997
$wsMetrics.Range("B4").Value  = 118564.67
$wsMetrics.Range("B5").Value  = 15693
$wsMetrics.Range("B6").Value  = 5590401.04
$wsMetrics.Range("B7").Value  = 4732270.6100000003
$wsMetrics.Range("B8").Value  = 1650521.5500000003
$wsMetrics.Range("B9").Value  = 218400
$wsMetrics.Range("B10").Value = 34055782.030000001
$wsMetrics.Range("B11").Value = 32007545.770000003
$wsMetrics.Range("B12").Value = 11932243.589999994
$wsMetrics.Range("B13").Value = 1316030

# Move the Metrics sheet's own selection to D14 (as in the saved file).
$wsMetrics.Activate()
$wsMetrics.Range("D14").Select()

# --- Restore the workbook's active sheet ("today") + its selection -----
$wsToday = $wb.Worksheets.Item("today")
$wsToday.Activate()
$wsToday.Range("D6").Select()
